# Update the marksheet's correct/total marks values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: number right answers -> 5
$ws.Range("B11").Value = 5

# "Total" row: total marks -> 100
$ws.Range("B12").Value = 100

# Correct/total marks summary text
$ws.Range("E12").Value = "100/140"
